# Implements SME feedback from Lisa Murphy:
#   " and agree with it. Write your initials where the forms tell you too."
# becomes
#   " and agree with it. Print your forms and write your initials where the
#    forms tell you to. If you make any changes after you print the forms,
#    you should also write your initials by what you change."
# and the (hidden) _GoBack bookmark, which originally sat right before the
# trailing "." of that sentence, ends up at the very end of the paragraph,
# after all of the newly added text.

$d = $word.ActiveDocument

# 1. Find and remove "Write your initials where the forms tell you too."
#    This leaves " and agree with it. " followed immediately by the
#    bookmark and then the lone "." run that used to close the sentence.
#    (Replace:=0/wdReplaceNone just positions $target on the match; we
#    delete it explicitly below.)
$target = $d.Content
$target.Find.ClearFormatting()
$found = $target.Find.Execute("Write your initials where the forms tell you too.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the sentence to replace."
}
$target.Delete()

# 2. The _GoBack bookmark is collapsed (Start = End) and now sits exactly
#    where that sentence used to be, immediately before the remaining lone
#    "." run. Remove that "." too -- it will be re-added as part of the new
#    combined sentence below.
$hasBookmark = $d.Bookmarks.Exists("_GoBack")
if ($hasBookmark) {
    $bookmark = $d.Bookmarks.Item("_GoBack")
    $period = $d.Range($bookmark.Start, $bookmark.Start + 1)
    if ($period.Text -eq ".") {
        $period.Delete()
    }
}

# 3. Insert the replacement text exactly at the (now shifted) bookmark
#    position. Text inserted right at a collapsed bookmark lands *before*
#    it and pushes the bookmark forward, so this both restores the
#    sentence and relocates _GoBack to sit after all of the new text --
#    matching the target layout.
$newText = "Print your forms and write your initials where the forms tell you to. If you make any changes after you print the forms, you should also write your initials by what you change."
if ($hasBookmark) {
    $bookmark = $d.Bookmarks.Item("_GoBack")
    $insertionPoint = $d.Range($bookmark.Start, $bookmark.Start)
    $insertionPoint.InsertAfter($newText)
} else {
    $target.Collapse(0)
    $target.InsertAfter($newText)
}

Write-Output "Updated Step 1 instructions and relocated _GoBack bookmark."
